# Apply country-data & header updates to the "Pais" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 16:52"

# --- Country name realignment (Republica Dominicana moved up before Bielorrusia) ---
$ws.Range("A48").Value = "Republica Dominicana"
$ws.Range("A49").Value = "Bielorrusia"
$ws.Range("A50").Value = "Singapur"

# --- Updated statistics values ---
$ws.Range("B4").Value = 588465
$ws.Range("C4").Value = 1524
$ws.Range("D4").Value = 37326
$ws.Range("E4").Value = 527428
$ws.Range("G4").Value = 71
$ws.Range("H4").Value = 23711
$ws.Range("B15").Value = 25913
$ws.Range("C15").Value = 225
$ws.Range("E15").Value = 11051
$ws.Range("G15").Value = 24
$ws.Range("H15").Value = 1162
$ws.Range("E31").Value = 5477
$ws.Range("G31").Value = 20
$ws.Range("H31").Value = 351
$ws.Range("B48").Value = 3286
$ws.Range("C48").Value = 119
$ws.Range("D48").Value = 162
$ws.Range("E48").Value = 2941
$ws.Range("F48").Value = 143
$ws.Range("G48").Value = 6
$ws.Range("H48").Value = 183
$ws.Range("B49").Value = 3281
$ws.Range("C49").Value = 362
$ws.Range("D49").Value = 203
$ws.Range("E49").Value = 3045
$ws.Range("F49").Value = 57
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 33
$ws.Range("B50").Value = 3252
$ws.Range("C50").Value = 334
$ws.Range("D50").Value = 611
$ws.Range("E50").Value = 2631
$ws.Range("F50").Value = 29
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 10
$ws.Range("B59").Value = 1934
$ws.Range("C59").Value = 222
$ws.Range("E59").Value = 1764
$ws.Range("B61").Value = 1720
$ws.Range("C61").Value = 9
$ws.Range("D61").Value = 989
$ws.Range("E61").Value = 723
$ws.Range("F61").Value = 8
$ws.Range("E74").Value = 943
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 26
$ws.Range("F85").Value = 29
